# checkpointing - working our hourly check for weather
#
# Insert a new "e022a Rain Roll" event row above the existing row 36
# ("e023 Call for Artillery Support"), pushing every subsequent row down
# by one. Then point the view at the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36 (existing rows 36.. shift down to 37..)
$ws.Rows.Item(36).Insert()

# New event id + rich description text for the inserted row.
$ws.Range("A36").Value = "e022a"
$ws.Range("B36").Value = "<Bold>e022a Rain Roll</Bold> `n<InlineUIContainer><Button Content='Weather' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    `n<LineBreak/><LineBreak/>"

# Match the row height used for this new entry.
$ws.Rows.Item(36).RowHeight = 45

# Move the viewport / selection onto the newly inserted row, like the author did.
$win = $excel.ActiveWindow
$win.ScrollRow = 34
$win.ScrollColumn = 1
[void]$ws.Range("B37").Select()
